$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WIP")

# Row 7: Campaign (group) changes from "COP 2026" to "COP 2027"
$ws.Range("A7").Value = "COP 2027"

# Row 12: Campaign (group) changes from "COP 2025" to "COP 2028"
$ws.Range("A12").Value = "COP 2028"

# Widen column C so the longer campaign labels fit
$ws.Columns.Item(3).ColumnWidth = 27.3

# Move/record the active selection to A14
$ws.Range("A14").Select()
